$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "zmiany nr 2"
$ws.Range("E10").Select()
